# Planets.xlsx: add an "ID" column, add "size"/"x"/"y" columns after "name",
# and append a small syntax-highlighted code snippet (rows 9-13, column C)
# showing how the planet constructor call looks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert the new "ID" column at the very left (A), and three new
#    columns ("size", "x", "y") right after the existing "name" column.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).Insert()
$ws.Range("C1:E1").EntireColumn.Insert()

# The insert operation copies the left neighbour's formatting onto the
# freshly inserted cells - strip that back off so the new cells fall
# back to the default (unstyled) look, matching the source data columns.
$ws.Range("A1:A6").ClearFormats()
$ws.Range("C1:E6").ClearFormats()

# ---------------------------------------------------------------------
# 2. Fill in the new "ID" column.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "ID"
$ws.Range("A2").Value = "Planet1"
$ws.Range("A3").Value = "Planet2"
$ws.Range("A4").Value = "Planet3"
$ws.Range("A5").Value = "Planet4"
$ws.Range("A6").Value = "Planet5"

# ---------------------------------------------------------------------
# 3. Fill in the new "size" / "x" / "y" columns.
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "size"
$ws.Range("D1").Value = "x"
$ws.Range("E1").Value = "y"

$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 650
$ws.Range("E2").Value = 90

$ws.Range("C3").Value = 80
$ws.Range("D3").Value = 350
$ws.Range("E3").Value = 190

$ws.Range("C4").Value = 150
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 40

$ws.Range("C5").Value = 70
$ws.Range("D5").Value = 500
$ws.Range("E5").Value = 300

$ws.Range("C6").Value = 200
$ws.Range("D6").Value = 810
$ws.Range("E6").Value = 300

# Column A ("ID") is a bit narrower than auto-fit would make it.
$ws.Columns.Item(1).ColumnWidth = 10

# ---------------------------------------------------------------------
# 4. Add a few rows below the table with a small syntax-highlighted
#    snippet (Consolas 7pt) illustrating the constructor call used to
#    build each Planet object from the row above: a teal "type" line,
#    a green "string" line, a grey "separator" line, repeated.
# ---------------------------------------------------------------------
$ws.Range("C9:C13").Font.Name = "Consolas"
$ws.Range("C9:C13").Font.Size = 7

$ws.Range("C12").Value = ", , "

$ws.Range("C9").Font.Color = 11062965   # FFB5CEA8 (string green)
$ws.Range("C10").Font.Color = 13948116  # FFD4D4D4 (plain grey)
$ws.Range("C11").Font.Color = 11585870  # FF4EC9B0 (type teal)

# Re-use the exact formatting already built above for the repeating rows
# instead of re-deriving new fonts for them.
$ws.Range("C10").Copy()
$ws.Range("C12").PasteSpecial(-4122)

$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5. Final cursor position, matching where editing left off.
# ---------------------------------------------------------------------
$ws.Range("F12").Select()

Write-Output "Planets.xlsx updated: ID/size/x/y columns + snippet rows added"
